$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-44 and 47-51: only Price (D) and Volume(1h) (E) columns change.
# Price values are stored as text in the sheet (e.g. "30.815.49", "1.000"),
# so a leading apostrophe forces Excel to keep them as text instead of
# reinterpreting number-looking strings (like "1.000") as numeric values.
$updates = @{
    2  = @("30.815.49", "  +1.32%  ")
    3  = @("1.887.30", "  +1.89%  ")
    4  = @("0.9999", "  +0.05%  ")
    5  = @("238.45", "  +2.02%  ")
    6  = @("1.000", "  +0.05%  ")
    7  = @("0.4775", "  +2.04%  ")
    8  = @("0.2876", "  +5.45%  ")
    9  = @("0.06574", "  +4.24%  ")
    10 = @("18.84", "  +15.09%  ")
    11 = @("97.66", "  +16.27%  ")
    12 = @("1.879.59", "  +1.88%  ")
    13 = @("0.07594", "  +1.88%  ")
    14 = @("5.113", "  +3.51%  ")
    15 = @("0.6580", "  +6.09%  ")
    16 = @("307.64", "  +33.77%  ")
    17 = @("30.810.21", "  +1.48%  ")
    18 = @("13.16", "  +6.17%  ")
    19 = @("1.000", "  -0.04%  ")
    20 = @("0.000007569", "  +3.71%  ")
    21 = @("2.107.57", "  +1.50%  ")
    22 = @("1.001", "  +0.19%  ")
    23 = @("5.122", "  +4.34%  ")
    24 = @("6.154", "  +4.55%  ")
    25 = @("9.276", "  +1.32%  ")
    26 = @("166.24", "  -0.52%  ")
    27 = @("20.27", "  +13.46%  ")
    28 = @("1.944", "  +3.81%  ")
    29 = @("0.1072", "  +5.01%  ")
    30 = @("1.354", "  -1.57%  ")
    31 = @("4.173", "  +2.20%  ")
    32 = @("3.972", "  +4.10%  ")
    33 = @("0.05033", "  +2.96%  ")
    34 = @("1.171", "  +2.77%  ")
    35 = @("0.7286", "  +3.65%  ")
    36 = @("2.710", "  +0.57%  ")
    37 = @("0.01943", "  +0.89%  ")
    38 = @("2.703", "  +1.40%  ")
    39 = @("2.078", "  +6.57%  ")
    40 = @("0.9030", "  +3.39%  ")
    41 = @("107.74", "  +2.05%  ")
    42 = @("1.000", "  +0.03%  ")
    43 = @("0.4200", "  +3.90%  ")
    44 = @("5.626", "  +1.91%  ")
    47 = @("8.994", "  +4.87%  ")
    48 = @("0.1222", "  +1.30%  ")
    49 = @("34.74", "  +3.99%  ")
    50 = @("0.05615", "  +1.25%  ")
    51 = @("1.386", "  +2.62%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = "'" + $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

# Rows 45 and 46 swap places (Aave <-> Aptos) with new values.
$ws.Cells.Item(45, 2).Value = "Aptos"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(45, 4).Value = "'7.349"
$ws.Cells.Item(45, 5).Value = "  +4.15%  "

$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).Value = "'65.59"
$ws.Cells.Item(46, 5).Value = "  +7.53%  "
